# Commit: "Update countries & provincias Spain"
#
# The "Pais" sheet is a COVID-19 dashboard whose rows are kept sorted by
# column B ("Casos totales") descending. This update refreshes the daily
# figures for several countries. Because a handful of the refreshed totals
# cross a neighbouring country total, the re-sort also shuffles the country
# name (col A) and stats (cols B:H) of the rows sitting in between, each down
# or up by one position (a couple of exactly-tied rows simply swap places).
# Below is the net per-cell effect of that refresh + re-sort, written as plain
# final values (row numbers / column letters match the saved worksheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: -> Alemania
$ws.Range("B8").Value = 135383
$ws.Range("C8").Value = 630
$ws.Range("D8").Value = 77000
$ws.Range("E8").Value = 54533
$ws.Range("F8").Value = 4288
$ws.Range("G8").Value = 46
$ws.Range("H8").Value = 3850

# Row 9: -> Reino Unido
$ws.Range("B9").Value = 103093
$ws.Range("C9").Value = 4617
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 89020
$ws.Range("F9").Value = 1559
$ws.Range("G9").Value = 861
$ws.Range("H9").Value = 13729

# Row 18: -> Suiza
$ws.Range("B18").Value = 26732
$ws.Range("C18").Value = 396
$ws.Range("D18").Value = 15400
$ws.Range("E18").Value = 10063
$ws.Range("F18").Value = 386
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = 1269

# Row 42: -> Serbia
$ws.Range("A42").Value = 'Serbia'
$ws.Range("B42").Value = 5318
$ws.Range("C42").Value = 445
$ws.Range("D42").Value = 443
$ws.Range("E42").Value = 4772
$ws.Range("F42").Value = 120
$ws.Range("G42").Value = 4
$ws.Range("H42").Value = 103

# Row 43: -> Malasia
$ws.Range("A43").Value = 'Malasia'
$ws.Range("B43").Value = 5182
$ws.Range("C43").Value = 110
$ws.Range("D43").Value = 2766
$ws.Range("E43").Value = 2332
$ws.Range("F43").Value = 56
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = 84

# Row 51: -> Finlandia
$ws.Range("B51").Value = 3369
$ws.Range("C51").Value = 132
$ws.Range("D51").Value = 1700
$ws.Range("E51").Value = 1594
$ws.Range("F51").Value = 76
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 75

# Row 102: -> San Marino
$ws.Range("A102").Value = 'San Marino'
$ws.Range("B102").Value = 426
$ws.Range("C102").Value = 33
$ws.Range("D102").Value = 55
$ws.Range("E102").Value = 333
$ws.Range("F102").Value = 15
$ws.Range("G102").Value = 2
$ws.Range("H102").Value = 38

# Row 103: -> Malta
$ws.Range("A103").Value = 'Malta'
$ws.Range("B103").Value = 412
$ws.Range("C103").Value = 13
$ws.Range("D103").Value = 82
$ws.Range("E103").Value = 327
$ws.Range("F103").Value = 4
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 3

# Row 104: -> Nigeria
$ws.Range("A104").Value = 'Nigeria'
$ws.Range("B104").Value = 407
$ws.Range("C104").Value = 0
$ws.Range("D104").Value = 128
$ws.Range("E104").Value = 267
$ws.Range("F104").Value = 2
$ws.Range("G104").Value = 0
$ws.Range("H104").Value = 12

# Row 105: -> Guinea
$ws.Range("A105").Value = 'Guinea'
$ws.Range("B105").Value = 404
$ws.Range("C105").Value = 0
$ws.Range("D105").Value = 31
$ws.Range("E105").Value = 372
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 1

# Row 106: -> Jordania
$ws.Range("A106").Value = 'Jordania'
$ws.Range("B106").Value = 401
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 250
$ws.Range("E106").Value = 144
$ws.Range("F106").Value = 5
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 7

# Row 107: -> Taiwan
$ws.Range("A107").Value = 'Taiwan'
$ws.Range("B107").Value = 395
$ws.Range("C107").Value = 0
$ws.Range("D107").Value = 155
$ws.Range("E107").Value = 234
$ws.Range("F107").Value = 0
$ws.Range("G107").Value = 0
$ws.Range("H107").Value = 6

# Row 142: -> Bermudas
$ws.Range("A142").Value = 'Bermudas'
$ws.Range("B142").Value = 81
$ws.Range("C142").Value = 0
$ws.Range("D142").Value = 33
$ws.Range("E142").Value = 43
$ws.Range("F142").Value = 3
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 5

# Row 143: -> Togo
$ws.Range("A143").Value = 'Togo'
$ws.Range("B143").Value = 81
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 35
$ws.Range("E143").Value = 43
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 3

# Row 182: -> Malaui
$ws.Range("A182").Value = 'Malaui'
$ws.Range("B182").Value = 16
$ws.Range("C182").Value = 0
$ws.Range("D182").Value = 0
$ws.Range("E182").Value = 14
$ws.Range("F182").Value = 1
$ws.Range("G182").Value = 0
$ws.Range("H182").Value = 2

# Row 183: -> Nepal
$ws.Range("A183").Value = 'Nepal'
$ws.Range("B183").Value = 16
$ws.Range("C183").Value = 0
$ws.Range("D183").Value = 2
$ws.Range("E183").Value = 14
$ws.Range("F183").Value = 0
$ws.Range("G183").Value = 0
$ws.Range("H183").Value = 0

# Row 184: -> Namibia
$ws.Range("B184").Value = 16
$ws.Range("C184").Value = 0
$ws.Range("D184").Value = 4
$ws.Range("E184").Value = 12
$ws.Range("F184").Value = 0
$ws.Range("G184").Value = 0
$ws.Range("H184").Value = 0

# Row 185: -> Suazilandia
$ws.Range("A185").Value = 'Suazilandia'
$ws.Range("B185").Value = 16
$ws.Range("C185").Value = 1
$ws.Range("D185").Value = 8
$ws.Range("E185").Value = 8
$ws.Range("F185").Value = 0
$ws.Range("G185").Value = 0
$ws.Range("H185").Value = 0

# Row 186: -> Dominica
$ws.Range("A186").Value = 'Dominica'
$ws.Range("B186").Value = 16
$ws.Range("C186").Value = 0
$ws.Range("D186").Value = 8
$ws.Range("E186").Value = 8
$ws.Range("F186").Value = 0
$ws.Range("G186").Value = 0
$ws.Range("H186").Value = 0

# Row 190: -> San Cristobal y Nieves
$ws.Range("A190").Value = 'San Cristobal y Nieves'
$ws.Range("B190").Value = 14
$ws.Range("C190").Value = 0
$ws.Range("D190").Value = 0
$ws.Range("E190").Value = 14
$ws.Range("F190").Value = 0
$ws.Range("G190").Value = 0
$ws.Range("H190").Value = 0

# Row 191: -> Granada
$ws.Range("A191").Value = 'Granada'
$ws.Range("B191").Value = 14
$ws.Range("C191").Value = 0
$ws.Range("D191").Value = 0
$ws.Range("E191").Value = 14
$ws.Range("F191").Value = 2
$ws.Range("G191").Value = 0
$ws.Range("H191").Value = 0

# Row 209: -> Sudan del Sur
$ws.Range("A209").Value = 'Sudan del Sur'
$ws.Range("B209").Value = 4
$ws.Range("C209").Value = 0
$ws.Range("D209").Value = 0
$ws.Range("E209").Value = 4
$ws.Range("F209").Value = 0
$ws.Range("G209").Value = 0
$ws.Range("H209").Value = 0

# Row 210: -> Santo Tome y Principe
$ws.Range("A210").Value = 'Santo Tome y Principe'
$ws.Range("B210").Value = 4
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 0
$ws.Range("E210").Value = 4
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 215: -> Yemen
$ws.Range("A215").Value = 'Yemen'
$ws.Range("B215").Value = 1
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 0
$ws.Range("E215").Value = 1
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

# Row 216: -> San Pedro y Miquelon
$ws.Range("A216").Value = 'San Pedro y Miquelon'
$ws.Range("B216").Value = 1
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 0
$ws.Range("E216").Value = 1
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 0
